$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.992.05"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "'2.988.96"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'559.47"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").Value = "'136.57"
$ws.Range("E6").Value = "  +3.58%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "'0.522"
$ws.Range("E8").Value = "  +1.52%  "

$ws.Range("D9").Value = "'2.979.77"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("E10").Value = "  +2.75%  "

$ws.Range("D11").Value = "'5.11"
$ws.Range("E11").Value = "  +6.25%  "

$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  +1.89%  "

$ws.Range("D13").Value = "'0.0000229"
$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("D14").Value = "'33.45"
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("E15").Value = "  +2.15%  "

$ws.Range("D16").Value = "'3.485.54"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("E17").Value = "  +6.82%  "

$ws.Range("D18").Value = "'2.992.04"
$ws.Range("E18").Value = "  +1.81%  "

$ws.Range("D19").Value = "'59.133.76"
$ws.Range("E19").Value = "  +3.00%  "

$ws.Range("D20").Value = "'427.39"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("D21").Value = "'13.60"
$ws.Range("E21").Value = "  +4.04%  "

$ws.Range("E22").Value = "  +5.34%  "

$ws.Range("D23").Value = "'7.08"
$ws.Range("E23").Value = "  +1.91%  "

$ws.Range("D24").Value = "'13.24"
$ws.Range("E24").Value = "  +2.15%  "

$ws.Range("D25").Value = "'80.36"
$ws.Range("E25").Value = "  +1.59%  "

$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  +0.30%  "

$ws.Range("D28").Value = "'2.16"
$ws.Range("E28").Value = "  +9.79%  "

$ws.Range("D29").Value = "'2.53"
$ws.Range("E29").Value = "  +1.92%  "

$ws.Range("D30").Value = "'7.75"
$ws.Range("E30").Value = "  +3.08%  "

$ws.Range("D31").Value = "'25.59"
$ws.Range("E31").Value = "  +2.05%  "

$ws.Range("D32").Value = "'5.99"
$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").Value = "'0.0986"
$ws.Range("E33").Value = "  -2.72%  "

$ws.Range("E34").Value = "  +6.10%  "

$ws.Range("D35").Value = "'5.90"
$ws.Range("E35").Value = "  +4.86%  "

$ws.Range("D36").Value = "0.0₃0749"
$ws.Range("E36").Value = "  +10.09%  "

$ws.Range("D37").Value = "'2.07"
$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").Value = "'48.86"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("D39").Value = "'8.65"
$ws.Range("E39").Value = "  +2.74%  "

$ws.Range("D40").Value = "'2.71"
$ws.Range("E40").Value = "  +6.13%  "

$ws.Range("D41").Value = "'396.62"
$ws.Range("E41").Value = "  +5.01%  "

$ws.Range("D42").Value = "'0.0350"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").Value = "'2.744.52"
$ws.Range("E43").Value = "  +3.73%  "

$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").Value = "'0.250"
$ws.Range("E45").Value = "  +4.44%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'35.09"
$ws.Range("E47").Value = "  +24.02%  "

$ws.Range("D48").Value = "'122.55"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("E49").Value = "  +0.84%  "

$ws.Range("D50").Value = "'1.98"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").Value = "'23.22"
$ws.Range("E51").Value = "  -0.34%  "
